# Strava activity sheet rework:
#  - Distance switched from meters to kilometers
#  - Moving Time split into a human "h:min" column and a raw-seconds column
#  - Average Speed switched from m/s to a pace in min/km
# This inserts one new column (for "Moving Time (seconds)") and rewrites the
# header row + the single data row to match the new layout.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Activity")

# Insert a brand-new column at D; everything from D onward (old Elapsed
# Time, Type, Start Date, Average Speed, Max Speed) slides one to the right.
$ws.Columns("D:D").Insert()

# --- Header row -----------------------------------------------------------
$ws.Range("B1").Value = "Distance (km)"
$ws.Range("C1").Value = "Moving Time (h:min)"
$ws.Range("D1").Value = "Moving Time (seconds)"
$ws.Range("H1").Value = "Average Speed (min/km)"

# Copy the style used by the rest of the header row onto the newly
# inserted D1 cell so it keeps the bold/centered/bordered look.
$ws.Range("C1").Copy()
$ws.Range("D1").PasteSpecial(-4122)  # xlPasteFormats
$excel.CutCopyMode = 0

# --- Data row ---------------------------------------------------------------
$ws.Range("B2").Value = 4.58
$ws.Range("C2").Value = "01:12:11"
$ws.Range("D2").Value = 4331
$ws.Range("H2").Value = "15:46"
$ws.Range("I2").Value = 11.71602432
